$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 28.634464
$ws.Range("H2").Value = 85.903392
$ws.Range("I2").Value = 0.3570833578723188
$ws.Range("J2").Value = 0.3570833578723188
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 112.513392
$ws.Range("N2").Value = 337.540176
$ws.Range("O2").Value = 0.3275312977368564
$ws.Range("P2").Value = 0.3275312977368564
$ws.Range("Q2").Value = 3221.760672741887
$ws.Range("R2").Value = 28995.84605467699
$ws.Range("S2").Value = 0.1169559756041549
$ws.Range("T2").Value = 0.1169559756041549

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 28.634464
$ws.Range("H3").Value = 85.903392
$ws.Range("I3").Value = 0.3570833578723188
$ws.Range("J3").Value = 0.3570833578723188
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 106.314466
$ws.Range("N3").Value = 318.943398
$ws.Range("O3").Value = 0.3094859589441663
$ws.Range("P3").Value = 0.3094859589441664
$ws.Range("Q3").Value = 3044.257749356223
$ws.Range("R3").Value = 27398.31974420602
$ws.Range("S3").Value = 0.1105122854341175
$ws.Range("T3").Value = 0.1105122854341175

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 28.634464
$ws.Range("H4").Value = 85.903392
$ws.Range("I4").Value = 0.3570833578723188
$ws.Range("J4").Value = 0.3570833578723188
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 124.6916553333333
$ws.Range("N4").Value = 374.074966
$ws.Range("O4").Value = 0.3629827433189773
$ws.Range("P4").Value = 0.3629827433189773
$ws.Range("Q4").Value = 3570.478715742741
$ws.Range("R4").Value = 32134.30844168467
$ws.Range("S4").Value = 0.1296150968340464
$ws.Range("T4").Value = 0.1296150968340464

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 36.74939233333333
$ws.Range("H5").Value = 110.248177
$ws.Range("I5").Value = 0.4582797992710433
$ws.Range("J5").Value = 0.4582797992710433
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 112.513392
$ws.Range("N5").Value = 337.540176
$ws.Range("O5").Value = 0.3275312977368564
$ws.Range("P5").Value = 0.3275312977368564
$ws.Range("Q5").Value = 4134.798785362128
$ws.Range("R5").Value = 37213.18906825915
$ws.Range("S5").Value = 0.1501009773818309
$ws.Range("T5").Value = 0.1501009773818309

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 36.74939233333333
$ws.Range("H6").Value = 110.248177
$ws.Range("I6").Value = 0.4582797992710433
$ws.Range("J6").Value = 0.4582797992710433
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 106.314466
$ws.Range("N6").Value = 318.943398
$ws.Range("O6").Value = 0.3094859589441663
$ws.Range("P6").Value = 0.3094859589441664
$ws.Range("Q6").Value = 3906.992021742827
$ws.Range("R6").Value = 35162.92819568545
$ws.Range("S6").Value = 0.1418311631421389
$ws.Range("T6").Value = 0.1418311631421389

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 36.74939233333333
$ws.Range("H7").Value = 110.248177
$ws.Range("I7").Value = 0.4582797992710433
$ws.Range("J7").Value = 0.4582797992710433
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 124.6916553333333
$ws.Range("N7").Value = 374.074966
$ws.Range("O7").Value = 0.3629827433189773
$ws.Range("P7").Value = 0.3629827433189773
$ws.Range("Q7").Value = 4582.342562537443
$ws.Range("R7").Value = 41241.08306283698
$ws.Range("S7").Value = 0.1663476587470736
$ws.Range("T7").Value = 0.1663476587470736

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 14.80600233333333
$ws.Range("H8").Value = 44.418007
$ws.Range("I8").Value = 0.1846368428566379
$ws.Range("J8").Value = 0.1846368428566379
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 112.513392
$ws.Range("N8").Value = 337.540176
$ws.Range("O8").Value = 0.3275312977368564
$ws.Range("P8").Value = 0.3275312977368564
$ws.Range("Q8").Value = 1665.873544483248
$ws.Range("R8").Value = 14992.86190034923
$ws.Range("S8").Value = 0.06047434475087062
$ws.Range("T8").Value = 0.06047434475087062

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 14.80600233333333
$ws.Range("H9").Value = 44.418007
$ws.Range("I9").Value = 0.1846368428566379
$ws.Range("J9").Value = 0.1846368428566379
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 106.314466
$ws.Range("N9").Value = 318.943398
$ws.Range("O9").Value = 0.3094859589441663
$ws.Range("P9").Value = 0.3094859589441664
$ws.Range("Q9").Value = 1574.092231663087
$ws.Range("R9").Value = 14166.83008496779
$ws.Range("S9").Value = 0.05714251036790991
$ws.Range("T9").Value = 0.05714251036790993

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 14.80600233333333
$ws.Range("H10").Value = 44.418007
$ws.Range("I10").Value = 0.1846368428566379
$ws.Range("J10").Value = 0.1846368428566379
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 124.6916553333333
$ws.Range("N10").Value = 374.074966
$ws.Range("O10").Value = 0.3629827433189773
$ws.Range("P10").Value = 0.3629827433189773
$ws.Range("Q10").Value = 1846.184939812529
$ws.Range("R10").Value = 16615.66445831276
$ws.Range("S10").Value = 0.06701998773785733
$ws.Range("T10").Value = 0.06701998773785735
